# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.316.90"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.875.72"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.7135"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'242.47"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.08038"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("D9").Value = "'0.3151"
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").Value = "'25.01"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "'0.08222"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "1.879.84"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "'94.89"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").Value = "'5.247"
$ws.Range("D15").Value = "'0.7123"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'6.410"
$ws.Range("E16").Value = "  +5.94%  "
$ws.Range("D17").Value = "'0.000008512"
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("D18").Value = "29.314.70"
$ws.Range("D19").Value = "'243.91"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.26"
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.128.32"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'7.762"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'0.1560"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'9.038"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'162.37"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'18.52"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "'1.505"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'4.407"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'4.306"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'0.05380"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  -8.57%  "
$ws.Range("D34").Value = "'1.938"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.7640"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "'1.179"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "'0.01876"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "1.261.75"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'2.755"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "'6.427"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'113.04"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.9111"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").Value = "'73.95"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("E45").Value = "  +9.55%  "
$ws.Range("D47").Value = "2.023.80"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'1.798"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'9.491"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "'0.4350"
$ws.Range("E51").Value = "  +1.10%  "
